$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Target values scraped from the commit diff (cell ref -> new text value).
$updates = [ordered]@{
    'D2' = '63.958.97'
    'D3' = '3.066.77'
    'D5' = '559.92'
    'D6' = '143.03'
    'D7' = '1.00'
    'D8' = '3.066.83'
    'D11' = '6.10'
    'D12' = '0.482'
    'D14' = '35.52'
    'D15' = '3.568.75'
    'D16' = '64.012.80'
    'D17' = '3.066.73'
    'D19' = '6.80'
    'D20' = '488.11'
    'D23' = '14.68'
    'D25' = '82.64'
    'D26' = '1.00'
    'D28' = '8.22'
    'D29' = '2.06'
    'D30' = '1.00'
    'D31' = '26.51'
    'D32' = '1.16'
    'D34' = '5.74'
    'D35' = '6.28'
    'D36' = '55.04'
    'D37' = '0.0413'
    'D38' = '445.73'
    'D39' = '0.0817'
    'D40' = '3.048.36'
    'D41' = '2.77'
    'D42' = '8.37'
    'D45' = '28.21'
    'D50' = '117.92'
    'D51' = '2.16'
    'E2' = '  -1.03%  '
    'E3' = '  -0.70%  '
    'E4' = '  +0.04%  '
    'E5' = '  -0.04%  '
    'E6' = '  -2.28%  '
    'E7' = '  -0.05%  '
    'E8' = '  -0.45%  '
    'E9' = '  +3.65%  '
    'E10' = '  +0.83%  '
    'E11' = '  -4.87%  '
    'E12' = '  +1.83%  '
    'E13' = '  +1.21%  '
    'E14' = '  +0.21%  '
    'E15' = '  -0.59%  '
    'E16' = '  -1.00%  '
    'E17' = '  -0.70%  '
    'E18' = '  -0.10%  '
    'E19' = '  +0.26%  '
    'E20' = '  +2.20%  '
    'E21' = '  +2.76%  '
    'E22' = '  +1.14%  '
    'E23' = '  +7.80%  '
    'E24' = '  -0.08%  '
    'E25' = '  +1.86%  '
    'E26' = '  -0.06%  '
    'E27' = '  +0.58%  '
    'E28' = '  +0.12%  '
    'E29' = '  -0.26%  '
    'E31' = '  +1.44%  '
    'E32' = '  +0.72%  '
    'E33' = '  +3.18%  '
    'E34' = '  +2.18%  '
    'E35' = '  +1.78%  '
    'E36' = '  +0.37%  '
    'E37' = '  +1.18%  '
    'E38' = '  -5.08%  '
    'E39' = '  -2.34%  '
    'E40' = '  +2.49%  '
    'E41' = '  -8.33%  '
    'E42' = '  +1.07%  '
    'E43' = '  +1.83%  '
    'E44' = '  +6.94%  '
    'E45' = '  -1.32%  '
    'E46' = '  +4.57%  '
    'E47' = '  -0.03%  '
    'E48' = '  +1.46%  '
    'E49' = '  -0.60%  '
    'E50' = '  -0.20%  '
    'E51' = '  +3.82%  '
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    if ($ref.Substring(0,1) -eq "D") {
        # Column D holds price text like "63.958.97" or "1.00" that must remain
        # literal text (not be reinterpreted as a number), so force a Text
        # number format before writing the value.
        $cell.NumberFormat = "@"
    }
    $cell.Value = $updates[$ref]
}
